{"js": "// Add two new list items (\"Notifica\u00e7\u00f5es\" and \"Emails de promo\u00e7\u00f5es/Atualiza\u00e7\u00f5es\n// de pre\u00e7os\") at the end of the document, right after the last existing\n// paragraph (\"Filtros\"). The new paragraphs inherit the list formatting\n// (style \"PargrafodaLista\", numId 1, ilvl 0) from the paragraph they are\n// inserted after, matching the rest of the bulleted/numbered list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document body (\"Filtros\") anchors the insert.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert \"Notifica\u00e7\u00f5es\" right after the last paragraph, then insert the\n// second new item right after that one \u2014 this keeps both new paragraphs in\n// the same list (same numId/ilvl) as \"Filtros\".\nconst notificacoesParagraph = lastParagraph.insertParagraph(\n  \"Notifica\u00e7\u00f5es\",\n  Word.InsertLocation.after\n);\nnotificacoesParagraph.insertParagraph(\n  \"Emails de promo\u00e7\u00f5es/Atualiza\u00e7\u00f5es de pre\u00e7os\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Add two new list items (\"Notifica\u00e7\u00f5es\" and \"Emails de promo\u00e7\u00f5es/Atualiza\u00e7\u00f5es\n# de pre\u00e7os\") at the end of the document, right after the last existing\n# paragraph (\"Filtros\"). InsertParagraphAfter() on the last paragraph's Range\n# creates a new paragraph that inherits the same list formatting (style\n# \"PargrafodaLista\", numId 1, ilvl 0) as the paragraph it follows.\n\n$d = $word.ActiveDocument\n\n# Insert the first new item (\"Notifica\u00e7\u00f5es\") right after the current last\n# paragraph (\"Filtros\").\n$last = $d.Paragraphs.Last\n$last.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Last\n$p1.Range.Text = \"Notifica\u00e7\u00f5es\"\n\n# Insert the second new item right after the one we just added.\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Last\n$p2.Range.Text = \"Emails de promo\u00e7\u00f5es/Atualiza\u00e7\u00f5es de pre\u00e7os\"\n"}
